$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy A13 into B13:D13 and A14 into B14:D14
$ws.Range("A13").Copy()
$ws.Range("B13:D13").PasteSpecial()

$ws.Range("A14").Copy()
$ws.Range("B14:D14").PasteSpecial()

$ws.Range("D14").Value = "pol_spec"

$ws.Range("D13:D14").Select()
